# Fix header labels on existing sheets
$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add a new "PO Forecast" sheet at the end of the workbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Cells.Item(1, 1).Value = "ds"
$wsForecast.Cells.Item(1, 2).Value = "PO_Forecast"
$wsForecast.Cells.Item(1, 3).Value = "yhat_lower"
$wsForecast.Cells.Item(1, 4).Value = "yhat_upper"

# Match the header style used on the other sheets (bold, thin border, centered/top)
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:B1").PasteSpecial(-4122)
$wsWeekly.Range("B1").Copy()
$wsForecast.Range("C1:D1").PasteSpecial(-4122)

# Forecast data rows
$data = @(
    @(44997.99999999999, 24, -122.5084982029286, 165.7560667117011),
    @(45004.99999999999, 27, -104.1014997462942, 164.6150147550311),
    @(45011.99999999999, 30, -108.9584949473099, 170.0206894828112),
    @(45018.99999999999, 33, -98.60091570531453, 178.8893984322968),
    @(45025.99999999999, 36, -103.5924026859547, 179.3124703629334),
    @(45032.99999999999, 39, -108.278176973358, 173.0480418608986),
    @(45039.99999999999, 42, -106.7896833745972, 176.4614542357628),
    @(45060.99999999999, 51, -95.53162470867485, 184.0516877825017),
    @(45081.99999999999, 61, -79.143932840007, 196.2043812206884),
    @(45088.99999999999, 64, -83.23082091882304, 197.9316033748386),
    @(45095.99999999999, 67, -68.58291814396028, 199.2769841776806),
    @(45109.99999999999, 73, -68.3468655568671, 219.9048900360181),
    @(45123.99999999999, 79, -65.30615139710221, 219.9919521727081),
    @(45137.99999999999, 85, -56.69021444719328, 222.1836427781349),
    @(45144.99999999999, 88, -55.58448381039601, 234.9345176061575),
    @(45151.99999999999, 91, -61.19064811898711, 231.2453721976895),
    @(45179.99999999999, 104, -45.07403408440374, 253.0406550463245),
    @(45186.99999999999, 107, -32.36675954738506, 242.8114297946497),
    @(45193.99999999999, 110, -25.35063189042251, 250.8975924458807),
    @(45200.99999999999, 113, -23.37044566613846, 253.3829497739792),
    @(45207.99999999999, 116, -17.26119893426268, 247.0524964879277),
    @(45214.99999999999, 119, -25.8075034265508, 251.9605912051957),
    @(45228.99999999999, 125, -14.25968297814287, 264.7635244934343),
    @(45235.99999999999, 128, -7.780005816248842, 275.5963393329988),
    @(45242.99999999999, 131, -11.11150688877862, 271.5341334867132),
    @(45249.99999999999, 134, -11.37290100411912, 275.1033219965252),
    @(45256.99999999999, 137, -5.751723508339137, 276.5316177119988),
    @(45263.99999999999, 140, -5.711528085765783, 270.002557991372),
    @(45270.99999999999, 144, -12.30921628001587, 280.6416854141876),
    @(45277.99999999999, 147, 0.2047854813700731, 295.6758390538133),
    @(45284.99999999999, 150, 9.576975918302468, 293.7236925213876),
    @(45291.99999999999, 153, 5.787236691105715, 296.1190686853122),
    @(45298.99999999999, 156, -1.213491491147797, 296.370373192424),
    @(45305.99999999999, 159, 15.27222542123963, 297.3173597436865),
    @(45312.99999999999, 162, 24.23588853259301, 297.5581836173918),
    @(45319.99999999999, 165, 24.0301702489393, 315.4760995047154),
    @(45326.99999999999, 168, 41.74235454205605, 310.3277627941716)
)

$r = 2
foreach ($row in $data) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# Dates in column A use the same custom date-time number format as the other sheets
$wsForecast.Range("A2:A38").NumberFormat = "YYYY-MM-DD HH:MM:SS"
